$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.24309515953064
$ws.Range("B1").Value = 2.774164915084839
$ws.Range("C1").Value = 2.933178901672363
$ws.Range("D1").Value = 2.496977090835571
$ws.Range("E1").Value = 0.8581372499465942
